$d = $word.ActiveDocument

# Locate the paragraph ending in "...or automate most of the jobs." (the
# last content paragraph before the trailing empty paragraph) so we can
# insert the new material right after it, ahead of that trailing blank
# paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*or automate most of the jobs.*") {
        $targetIndex = $i
    }
}

# Collapse to the *start* of the paragraph immediately following the target
# (rather than collapsing the target's own range to its end) -- this keeps
# the insertion point logically "before" the following paragraph mark so
# InsertXML splices new paragraphs in without clobbering the preceding run.
if ($targetIndex -ge 1 -and $targetIndex -lt $d.Paragraphs.Count) {
    $r = $d.Paragraphs($targetIndex + 1).Range
    $r.Collapse(1)
} else {
    $r = $d.Content
    $r.Collapse(0)
}

$insertXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t xml:space="preserve">big data can basically be classified into two categories, namely, data from the physical world, which is usually obtained through sensors, scientific </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t>experiments</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t xml:space="preserve"> and observations (such as biological data, neural data, astronomical data, and remote sensing data), and data from the human society, which is often acquired from such sources or domains as social networks, Internet, health, finance, economics, and transportation.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t>Actually, the</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t xml:space="preserve"> real challenges </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t>center</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t xml:space="preserve"> around the diversified data types (Variety), timely response requirements (Velocity), and uncertainties in the data (Veracity).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t xml:space="preserve">emi-structured or unstructured data </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_Hlk145255274"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t>(including text, images, video, and voice)</w:t>
      </w:r>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Timely responses are also challenging because there may not be enough resources to collect, store, and process the big data within a reasonable amount of time. Finally, distinguishing between true and false or reliable and unreliable data is especially challenging, even for the best data cleaning methods to eliminate some inherent unpredictability of data.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>
          <w:color w:val="2E2E2E"/>
        </w:rPr>
        <w:t>From the socio-economic point of view, big data is the core connotation and critical support of the so-called second economy, a concept proposed by the American economist W.B. Arthur in 2011 [12], which refers to the economic activities running on processor, connectors, sensors, and executors. It is estimated that by 2030 the size of the second economy will approach that of the first economy (namely, the traditional physical economy). </w:t>
      </w:r>
    </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($insertXml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
